# Automatic update [2026-02-23 23:49]: refresh DATA_EXTRACCIO timestamps
# and the handful of observation values (humitat/pressio/radiacio/temperatura)
# that shifted between the 23:18 and 23:48 meteo.cat extraction runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-23 23:48:13'
$ws.Range('E3').Value = '2026-02-23 23:48:15'
$ws.Range('H3').Value = '''33%'
$ws.Range('E4').Value = '2026-02-23 23:48:18'
$ws.Range('O4').Value = '11.7 °C'
$ws.Range('E5').Value = '2026-02-23 23:48:20'
$ws.Range('K5').Value = '16.3 MJ/m2'
$ws.Range('E6').Value = '2026-02-23 23:48:22'
$ws.Range('E7').Value = '2026-02-23 23:48:25'
$ws.Range('E8').Value = '2026-02-23 23:48:27'
$ws.Range('E9').Value = '2026-02-23 23:48:29'
$ws.Range('E10').Value = '2026-02-23 23:48:32'
$ws.Range('O10').Value = '10.3 °C'
$ws.Range('E11').Value = '2026-02-23 23:48:34'
$ws.Range('O11').Value = '8.4 °C'
$ws.Range('E12').Value = '2026-02-23 23:48:36'
$ws.Range('E13').Value = '2026-02-23 23:48:39'
$ws.Range('H13').Value = '''61%'
$ws.Range('J13').Value = '1027.1 hPa'
$ws.Range('O13').Value = '6.6 °C'
$ws.Range('E14').Value = '2026-02-23 23:48:41'
$ws.Range('H14').Value = '''78%'
$ws.Range('E15').Value = '2026-02-23 23:48:43'
$ws.Range('O15').Value = '12.0 °C'
$ws.Range('E16').Value = '2026-02-23 23:48:46'
$ws.Range('H16').Value = '''24%'
$ws.Range('E17').Value = '2026-02-23 23:48:48'
$ws.Range('H17').Value = '''41%'
$ws.Range('E18').Value = '2026-02-23 23:48:50'
$ws.Range('H18').Value = '''76%'
$ws.Range('J18').Value = '1024.8 hPa'
$ws.Range('O18').Value = '10.4 °C'
$ws.Range('E19').Value = '2026-02-23 23:48:53'
$ws.Range('H19').Value = '''49%'
$ws.Range('O19').Value = '12.1 °C'
$ws.Range('E20').Value = '2026-02-23 23:48:55'
$ws.Range('H20').Value = '''39%'
$ws.Range('E21').Value = '2026-02-23 23:48:57'
$ws.Range('E22').Value = '2026-02-23 23:48:59'
$ws.Range('E23').Value = '2026-02-23 23:49:02'
$ws.Range('E24').Value = '2026-02-23 23:49:04'
$ws.Range('E25').Value = '2026-02-23 23:49:06'
$ws.Range('H25').Value = '''29%'
$ws.Range('E26').Value = '2026-02-23 23:49:09'
$ws.Range('E27').Value = '2026-02-23 23:49:11'
$ws.Range('E28').Value = '2026-02-23 23:49:14'
$ws.Range('O28').Value = '10.7 °C'
$ws.Range('E29').Value = '2026-02-23 23:49:16'
$ws.Range('O29').Value = '10.3 °C'
$ws.Range('E30').Value = '2026-02-23 23:49:18'
$ws.Range('H30').Value = '''72%'
$ws.Range('J30').Value = '1024.4 hPa'
$ws.Range('K30').Value = '15.2 MJ/m2'
$ws.Range('E31').Value = '2026-02-23 23:49:21'
$ws.Range('J31').Value = '1023.6 hPa'
$ws.Range('K31').Value = '15.1 MJ/m2'
$ws.Range('E32').Value = '2026-02-23 23:49:23'
$ws.Range('N32').Value = '-2.0 °C 23:26 TU'
$ws.Range('O32').Value = '6.6 °C'
$ws.Range('E33').Value = '2026-02-23 23:49:25'
$ws.Range('J33').Value = '1025.5 hPa'
$ws.Range('O33').Value = '8.4 °C'
$ws.Range('E34').Value = '2026-02-23 23:49:27'
$ws.Range('H34').Value = '''45%'
$ws.Range('O34').Value = '3.8 °C'
$ws.Range('E35').Value = '2026-02-23 23:49:30'
$ws.Range('E36').Value = '2026-02-23 23:49:32'
$ws.Range('H36').Value = '''74%'
$ws.Range('O36').Value = '12.7 °C'
$ws.Range('E37').Value = '2026-02-23 23:49:35'
$ws.Range('O37').Value = '8.6 °C'
$ws.Range('E38').Value = '2026-02-23 23:49:37'
$ws.Range('O38').Value = '11.9 °C'
$ws.Range('E39').Value = '2026-02-23 23:49:39'
$ws.Range('H39').Value = '''28%'
$ws.Range('E40').Value = '2026-02-23 23:49:42'
$ws.Range('J40').Value = '1026.4 hPa'
$ws.Range('O40').Value = '8.2 °C'
$ws.Range('E41').Value = '2026-02-23 23:49:44'
$ws.Range('E42').Value = '2026-02-23 23:49:47'
$ws.Range('H42').Value = '''81%'
$ws.Range('O42').Value = '11.5 °C'
$ws.Range('E43').Value = '2026-02-23 23:49:49'
$ws.Range('O43').Value = '10.2 °C'
$ws.Range('E44').Value = '2026-02-23 23:49:51'
$ws.Range('N44').Value = '-0.7 °C 23:01 TU'
$ws.Range('O44').Value = '2.9 °C'
$ws.Range('E45').Value = '2026-02-23 23:49:54'
$ws.Range('E46').Value = '2026-02-23 23:49:56'
$ws.Range('H46').Value = '''75%'
$ws.Range('O46').Value = '9.8 °C'
